$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: add Work Done (C18) and Bugs/Limitations (E18) text for day 14
$ws.Range("C18").Value = "keyword matching using nltk. created a log for keywords found with article link and company name. pip freeze clean requirements.txt file."
$ws.Range("E18").Value = "Try converting .py to .exe"

# Row 18 grows taller to fit the wrapped text that was just added
$ws.Rows.Item(18).RowHeight = 45

# Row 19: add the date for day 15
$ws.Range("B19").Value = 43649
$ws.Range("B19").NumberFormat = "d-mmm-yy"

# Update the window scroll position / selection to match where the user ended up
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
$ws.Range("E19").Select()
